# LearningTree.xlsx — "Latest code for provide score for all Activities"
#
# The automation test data sheet (STAGE) stores one randomly generated
# "CourseDesigner#####" / "LearningCourse#####" pair per run in cells
# K2 / L2. Each re-run of the data generator advances those two values.
# Apply the latest generated pair to the STAGE sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAGE")

$ws.Range("K2").Value = "CourseDesigner70039"
$ws.Range("L2").Value = "LearningCourse45781"
